# Fixing how the timezone appears in the ISO timestamps
#
# The Timestamp_nominal_HST_ISO8601 (col F) and Timestamp_bottle_closure_GMT_ISO8601
# (col G) columns store ISO8601 timestamps whose trailing UTC-offset was written
# with a colon (e.g. "-10:00" / "-00:00"). This rewrites every such timestamp in
# rows 5-52 so the offset has no colon (e.g. "-1000" / "-0000"), matching the
# corrected ISO8601 representation used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-TimezoneColon([string]$cellRef) {
    $cell = $ws.Range($cellRef)
    $val = $cell.Value2
    if ($val -ne $null -and $val -match '^(\d{4}-\d{2}-\d{2}T\d{2}:\d{2}:\d{2}-)(\d{2}):(\d{2})$') {
        $newVal = $Matches[1] + $Matches[2] + $Matches[3]
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Rows 5-52 hold the per-sample timestamps in columns F and G.
Fix-TimezoneColon "F5"
Fix-TimezoneColon "F7"
Fix-TimezoneColon "G5"
Fix-TimezoneColon "G7"
Fix-TimezoneColon "F6"
Fix-TimezoneColon "F8"
Fix-TimezoneColon "G6"
Fix-TimezoneColon "G8"

for ($r = 9; $r -le 52; $r += 2) {
    Fix-TimezoneColon ("F" + $r)
    Fix-TimezoneColon ("G" + $r)
    Fix-TimezoneColon ("F" + ($r + 1))
    Fix-TimezoneColon ("G" + ($r + 1))
}

# Restore the sheet's on-save selection state (scrolled back to top, columns H:K
# selected) rather than leaving it parked on the last-edited cell.
$ws.Range("H1:K1048576").Select()
